$d = $word.ActiveDocument

# --- Insert the four new paragraphs right before the trailing bookmark
#     paragraph ("_GoBack"), which currently only holds the bookmark.
#     A trailing paragraph-break ("`r") is included so the final inserted
#     run ("Eine computerstimme ...") is treated as a brand new paragraph
#     (and therefore correctly inherits the sz/szCs run formatting from
#     the paragraph mark) instead of being merged into the formatting-less
#     leftover run of the original (empty) bookmark paragraph.
$bm = $d.Bookmarks.Item("_GoBack")
$insertPoint = $d.Range($bm.Start, $bm.Start)

$newText = "only_input_alias" + "`r" + `
    "Code wird nur weiter ausgeführt wenn der Input der von dem User erfragt wird dem Alias entspricht" + "`r" + `
    "computer_voice_alias" + "`r" + `
    "Eine computerstimme Stimme gibt das aus was im Alias ist " + "`r"

$insertPoint.InsertBefore($newText)

# --- The insert above leaves one extra empty paragraph (still carrying
#     the original "_GoBack" bookmark) after the new "Eine computerstimme
#     ..." paragraph. Remove that extra paragraph mark so the bookmark
#     becomes part of the same paragraph as the new text again, exactly
#     like it originally sat in the (now replaced) empty paragraph.
$mergeParagraph = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$pilcrow = $d.Range($mergeParagraph.Range.End - 1, $mergeParagraph.Range.End)
$pilcrow.Delete()

# --- Highlight the two new alias names exactly like the other *_alias
#     entries already in the document (dark yellow highlight).
$find1 = $d.Content
$find1.Find.ClearFormatting()
$find1.Find.Text = "only_input_alias"
$find1.Find.Execute() | Out-Null
$find1.HighlightColorIndex = 14

$find2 = $d.Content
$find2.Find.ClearFormatting()
$find2.Find.Text = "computer_voice_alias"
$find2.Find.Execute() | Out-Null
$find2.HighlightColorIndex = 14
